$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 235, shifting existing rows 235-251 down to 236-252.
$ws.Rows.Item(235).Insert()

# Copy style (date format) from the cell above (D234) into the new D235 cell.
$ws.Range("D234").Copy()
$ws.Range("D235").PasteSpecial(-4122) # xlPasteFormats

# Populate the new row 235 with the weekly data point.
$ws.Cells.Item(235, 1).Value = 3
$ws.Cells.Item(235, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(235, 3).Value = "Coquimbo"
$ws.Cells.Item(235, 4).Value = 44585
$ws.Cells.Item(235, 5).Value = 5
$ws.Cells.Item(235, 6).Value = 100112039
$ws.Cells.Item(235, 7).Value = "Ciboulette"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 160
$ws.Cells.Item(235, 11).Value = 1500
$ws.Cells.Item(235, 12).Value = 1500
$ws.Cells.Item(235, 13).Value = 1500
$ws.Cells.Item(235, 14).Value = "$/docena de atados"
$ws.Cells.Item(235, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(235, 16).Value = 500
$ws.Cells.Item(235, 17).Value = 3
$ws.Cells.Item(235, 18).Value = "Hortaliza"
